$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.388106333333333
$ws.Range("H2").Value = 28.164319
$ws.Range("I2").Value = 0.2414596449149976
$ws.Range("J2").Value = 0.2414596449149975
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 454.6116294582945
$ws.Range("R2").Value = 4091.504665124651
$ws.Range("S2").Value = 0.1524035938700456
$ws.Range("T2").Value = 0.1524035938700455

$ws.Range("G3").Value = 9.388106333333333
$ws.Range("H3").Value = 28.164319
$ws.Range("I3").Value = 0.2414596449149976
$ws.Range("J3").Value = 0.2414596449149975
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("N3").Value = 20.549744
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 64.30772726492623
$ws.Range("R3").Value = 578.769545384336
$ws.Range("S3").Value = 0.02155846466239284
$ws.Range("T3").Value = 0.02155846466239284

$ws.Range("G4").Value = 9.388106333333333
$ws.Range("H4").Value = 28.164319
$ws.Range("I4").Value = 0.2414596449149976
$ws.Range("J4").Value = 0.2414596449149975
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 201.3416281773673
$ws.Range("R4").Value = 1812.074653596306
$ws.Range("S4").Value = 0.06749758638255915
$ws.Range("T4").Value = 0.06749758638255914

$ws.Range("I5").Value = 0.5770971896641285
$ws.Range("J5").Value = 0.5770971896641284
$ws.Range("M5").Value = 48.42420966666666
$ws.Range("N5").Value = 145.272629
$ws.Range("O5").Value = 0.6311762527593259
$ws.Range("P5").Value = 0.6311762527593258
$ws.Range("Q5").Value = 1086.538058321797
$ws.Range("R5").Value = 9778.842524896168
$ws.Range("S5").Value = 0.3642500416501426
$ws.Range("T5").Value = 0.3642500416501425

$ws.Range("I6").Value = 0.5770971896641285
$ws.Range("J6").Value = 0.5770971896641284
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("O6").Value = 0.08928392431779728
$ws.Range("P6").Value = 0.08928392431779726
$ws.Range("Q6").Value = 153.6977687983467
$ws.Range("S6").Value = 0.05152550180598555
$ws.Range("T6").Value = 0.05152550180598553

$ws.Range("I7").Value = 0.5770971896641285
$ws.Range("J7").Value = 0.5770971896641284
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2795398229228769
$ws.Range("P7").Value = 0.2795398229228769
$ws.Range("Q7").Value = 481.2136944227799
$ws.Range("R7").Value = 4330.923249805019
$ws.Range("S7").Value = 0.1613216462080004
$ws.Range("T7").Value = 0.1613216462080004

$ws.Range("I8").Value = 0.181443165420874
$ws.Range("J8").Value = 0.1814431654208739
$ws.Range("M8").Value = 48.42420966666666
$ws.Range("N8").Value = 145.272629
$ws.Range("O8").Value = 0.6311762527593259
$ws.Range("P8").Value = 0.6311762527593258
$ws.Range("Q8").Value = 341.6147369681277
$ws.Range("R8").Value = 3074.532632713149
$ws.Range("S8").Value = 0.1145226172391377
$ws.Range("T8").Value = 0.1145226172391377

$ws.Range("I9").Value = 0.181443165420874
$ws.Range("J9").Value = 0.1814431654208739
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("O9").Value = 0.08928392431779728
$ws.Range("P9").Value = 0.08928392431779726
$ws.Range("Q9").Value = 48.32359295516267
$ws.Range("R9").Value = 434.912336596464
$ws.Range("S9").Value = 0.01619995784941888
$ws.Range("T9").Value = 0.01619995784941888

$ws.Range("I10").Value = 0.181443165420874
$ws.Range("J10").Value = 0.1814431654208739
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2795398229228769
$ws.Range("P10").Value = 0.2795398229228769
$ws.Range("S10").Value = 0.05072059033231736
$ws.Range("T10").Value = 0.05072059033231736
